# Atualiza datasets das ligas
# Adds the new team "Pepe Leal FC" to the league workbook:
#  - "Geral" and "Turno 2": append a new row (19) for "Texas Club 2026"
#    (the alphabetical list grows by one row because "Pepe Leal FC" now
#    takes the slot "Texas Club 2026" used to occupy).
#  - "Classif Turno 2" and each "Mês - ..." sheet: re-sequence the
#    standings column with "Pepe Leal FC" inserted, and append a new
#    row (19) for "Texas Club 2026".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheets "Geral" and "Turno 2": alphabetical team list, columns B.. are
# round-by-round scores. "Pepe Leal FC" is inserted alphabetically
# (between "NHU PORÃ SAF." and "Pontaç0 F.C."), so every row from there
# on shifts down by one, and a new last row (19) appears for the team
# that used to be last ("Texas Club 2026").
# ---------------------------------------------------------------
$alphaSheets = @(
    @{ Name = "Geral";    LastCol = "AM" },
    @{ Name = "Turno 2";  LastCol = "T" }
)

$alphaOrder = @(
    "bugredasmissões",
    "C R Juvenal",
    "Doug Leal F.C",
    "Esquadrão Gazembrino",
    "FBC Colorado",
    "GaúchoDaFronteira F.C",
    "GE Bebum",
    "Grêmio_Campeão_LA_27",
    "JV5 Tricolor Gaúcho",
    "La Primeira Patada Es Nuestra",
    "lsauer fc",
    "Medonho´s F.C.",
    "NHU PORÃ SAF.",
    "Pepe Leal FC",
    "Pontaç0 F.C.",
    "SC 100 Sono",
    "SC ÉoINTER!",
    "Texas Club 2026"
)

foreach ($s in $alphaSheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    # Copy the whole row 18 formatting down into row 19 first so the new
    # row picks up the same styling (bold/centered/bordered team-name cell).
    $ws.Range("A18:" + $s.LastCol + "18").Copy()
    $ws.Range("A19:" + $s.LastCol + "19").PasteSpecial(-4122)

    for ($i = 0; $i -lt $alphaOrder.Length; $i++) {
        $row = $i + 2
        $ws.Range("A" + $row).Value = $alphaOrder[$i]
        $ws.Range("B" + $row + ":" + $s.LastCol + $row).Value = 0
    }
}

# ---------------------------------------------------------------
# Sheets "Classif Turno 2" and the monthly sheets: single score column
# (B), team names in column A in a fixed standings order. Re-write the
# column A order to insert "Pepe Leal FC" and append the new row 19.
# ---------------------------------------------------------------
$standingsSheets = @(
    "Classif Turno 2",
    "Mês - Janeiro",
    "Mês - Fevereiro",
    "Mês - Março",
    "Mês - Abril",
    "Mês - Maio",
    "Mês - Julho"
)

$newOrder = @(
    "bugredasmissões",
    "C R Juvenal",
    "SC ÉoINTER!",
    "SC 100 Sono",
    "Pontaç0 F.C.",
    "Pepe Leal FC",
    "NHU PORÃ SAF.",
    "Medonho´s F.C.",
    "lsauer fc",
    "La Primeira Patada Es Nuestra",
    "JV5 Tricolor Gaúcho",
    "Grêmio_Campeão_LA_27",
    "GE Bebum",
    "GaúchoDaFronteira F.C",
    "FBC Colorado",
    "Esquadrão Gazembrino",
    "Doug Leal F.C",
    "Texas Club 2026"
)

foreach ($name in $standingsSheets) {
    $ws = $wb.Worksheets.Item($name)

    # Extend formatting (team-name style + border) from row 18 to the new
    # row 19 before writing values.
    $ws.Range("A18:B18").Copy()
    $ws.Range("A19:B19").PasteSpecial(-4122)

    for ($i = 0; $i -lt $newOrder.Length; $i++) {
        $row = $i + 2
        $ws.Range("A" + $row).Value = $newOrder[$i]
        $ws.Range("B" + $row).Value = 0
    }
}
